# edit.ps1
# Applies crypto price/volume/name/link updates to Sheet1 as described by the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    ,@('D', 2, '65.254.30')
    ,@('E', 2, '  -2.46%  ')
    ,@('D', 3, '3.382.66')
    ,@('E', 3, '  -2.58%  ')
    ,@('D', 4, '0.998')
    ,@('E', 4, '  -0.19%  ')
    ,@('D', 5, '592.48')
    ,@('E', 5, '  -1.92%  ')
    ,@('D', 6, '140.74')
    ,@('E', 6, '  -5.15%  ')
    ,@('D', 7, '0.998')
    ,@('E', 7, '  -0.31%  ')
    ,@('D', 8, '3.381.76')
    ,@('E', 8, '  -2.55%  ')
    ,@('D', 9, '0.466')
    ,@('E', 9, '  -3.60%  ')
    ,@('B', 10, 'Toncoin')
    ,@('C', 10, 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton')
    ,@('D', 10, '7.90')
    ,@('E', 10, '  +4.39%  ')
    ,@('B', 11, 'Dogecoin')
    ,@('C', 11, 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge')
    ,@('D', 11, '0.133')
    ,@('E', 11, '  -6.91%  ')
    ,@('D', 12, '0.405')
    ,@('E', 12, '  -4.80%  ')
    ,@('D', 13, '3.946.81')
    ,@('E', 13, '  -2.81%  ')
    ,@('B', 14, 'ShibaInu')
    ,@('C', 14, 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib')
    ,@('D', 14, '0.0000198')
    ,@('E', 14, '  -7.89%  ')
    ,@('B', 15, 'Avalanche')
    ,@('C', 15, 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax')
    ,@('D', 15, '29.55')
    ,@('E', 15, '  -7.25%  ')
    ,@('E', 16, '  -0.64%  ')
    ,@('D', 17, '65.221.62')
    ,@('E', 17, '  -2.54%  ')
    ,@('D', 18, '3.376.91')
    ,@('E', 18, '  -2.78%  ')
    ,@('D', 19, '10.32')
    ,@('E', 19, '  +2.11%  ')
    ,@('D', 20, '6.08')
    ,@('E', 20, '  -6.04%  ')
    ,@('D', 21, '14.61')
    ,@('E', 21, '  -5.46%  ')
    ,@('D', 22, '413.50')
    ,@('E', 22, '  -6.25%  ')
    ,@('D', 23, '0.576')
    ,@('E', 23, '  -6.12%  ')
    ,@('D', 24, '76.99')
    ,@('E', 24, '  -2.54%  ')
    ,@('E', 25, '  -0.08%  ')
    ,@('D', 26, '3.507.96')
    ,@('E', 26, '  -2.76%  ')
    ,@('D', 27, '0.0000107')
    ,@('E', 27, '  -11.89%  ')
    ,@('D', 28, '9.17')
    ,@('E', 28, '  -6.32%  ')
    ,@('D', 29, '7.71')
    ,@('E', 29, '  -8.70%  ')
    ,@('D', 30, '2.40')
    ,@('E', 30, '  -3.42%  ')
    ,@('E', 31, '  +0.01%  ')
    ,@('E', 32, '  -4.55%  ')
    ,@('E', 33, '  -10.13%  ')
    ,@('B', 34, 'EthereumClassic')
    ,@('C', 34, 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc')
    ,@('D', 34, '24.26')
    ,@('E', 34, '  -4.69%  ')
    ,@('B', 35, 'RenzoRestakedETH')
    ,@('C', 35, 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth')
    ,@('D', 35, '3.371.37')
    ,@('E', 35, '  -2.68%  ')
    ,@('E', 36, '  -0.09%  ')
    ,@('E', 37, '  -7.67%  ')
    ,@('D', 38, '5.48')
    ,@('E', 38, '  -10.00%  ')
    ,@('D', 39, '7.49')
    ,@('E', 39, '  -5.72%  ')
    ,@('D', 40, '0.999')
    ,@('E', 40, '  -0.04%  ')
    ,@('B', 41, 'Hedera')
    ,@('C', 41, 'https://coinranking.com/coin/jad286TjB+hedera-hbar')
    ,@('D', 41, '0.0851')
    ,@('E', 41, '  -4.70%  ')
    ,@('B', 42, 'Monero')
    ,@('C', 42, 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr')
    ,@('D', 42, '164.55')
    ,@('E', 42, '  -5.54%  ')
    ,@('D', 43, '0.867')
    ,@('E', 43, '  -2.32%  ')
    ,@('D', 44, '4.99')
    ,@('E', 44, '  -7.99%  ')
    ,@('D', 45, '1.91')
    ,@('E', 45, '  -10.90%  ')
    ,@('D', 46, '45.22')
    ,@('E', 46, '  -2.14%  ')
    ,@('D', 47, '26.38')
    ,@('E', 47, '  -10.46%  ')
    ,@('E', 48, '  -7.02%  ')
    ,@('D', 49, '7.00')
    ,@('E', 49, '  -6.52%  ')
    ,@('E', 50, '  -9.58%  ')
    ,@('B', 51, 'SuiNetwork')
    ,@('C', 51, 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui')
    ,@('D', 51, '0.910')
    ,@('E', 51, '  -8.18%  ')
)

foreach ($change in $changes) {
    $col = $change[0]
    $row = $change[1]
    $val = $change[2]
    $cell = $ws.Range("$col$row")
    # Preserve the cell's existing style/number format (General) and write
    # the value with a leading apostrophe so Excel stores it verbatim as
    # text instead of auto-converting look-alike numbers (e.g. "65.254.30",
    # "0.998", "7.90", "0.0000198") which would otherwise lose formatting
    # such as trailing/leading zeros and literal dot-grouping.
    $origStyle = $cell.Style
    $cell.Value = "'" + $val
    $cell.Style = $origStyle
}
